$wb = $excel.ActiveWorkbook

# --- DBD sheet (sheet1): field "D20"/"D22" type changed from DATE to TIMESTAMP,
#     and their length ("E20"/"E22") no longer applies, so it's cleared. ---
$wsDBD = $wb.Worksheets.Item("DBD")

$wsDBD.Range("D20").Value = "TIMESTAMP"
$wsDBD.Range("E20").ClearContents()

$wsDBD.Range("D22").Value = "TIMESTAMP"
$wsDBD.Range("E22").ClearContents()

# --- Selection / active sheet ends up on DBD, at the last-touched cell. ---
$wsDBD.Select()
$wsDBD.Range("E22").Select()
